$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 344607
$ws.Range("B1").Value = 281842
$ws.Range("A2").Value = 325140
$ws.Range("B2").Value = 291722
$ws.Range("A3").Value = 323035
$ws.Range("B3").Value = 255341
$ws.Range("A4").Value = 351548
$ws.Range("B4").Value = 301403

$ws.Range("I27").Select()
